# Auto-generated Excel COM-interop script applying the Tiamat_Profits diff.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets("ALC")
$ws_ARM = $wb.Worksheets("ARM")
$ws_CRP = $wb.Worksheets("CRP")
$ws_CUL = $wb.Worksheets("CUL")
$ws_GSM = $wb.Worksheets("GSM")
$ws_LTW = $wb.Worksheets("LTW")

# ALC row 19
$ws_ALC.Range("H19").Value = 775.7083
$ws_ALC.Range("I19").Value = 833.2222
$ws_ALC.Range("J19").Value = 741.2
$ws_ALC.Range("K19").Value = 833.2222
$ws_ALC.Range("L19").Value = 741.2
$ws_ALC.Range("M19").Value = -658.2222
$ws_ALC.Range("N19").Value = -1091.2

# ALC row 32
$ws_ALC.Range("H32").Value = 1156.6666
$ws_ALC.Range("I32").Value = 1485
$ws_ALC.Range("J32").Value = 500
$ws_ALC.Range("K32").Value = 1485
$ws_ALC.Range("L32").Value = 500
$ws_ALC.Range("M32").Value = -1159
$ws_ALC.Range("N32").Value = -1152

# ALC row 53
$ws_ALC.Range("H53").Value = 309.7143
$ws_ALC.Range("I53").Value = 26.5
$ws_ALC.Range("J53").Value = 687.3333
$ws_ALC.Range("K53").Value = 26.5
$ws_ALC.Range("L53").Value = 687.3333
$ws_ALC.Range("M53").Value = 610.5
$ws_ALC.Range("N53").Value = -1961.3333

# ALC row 55
$ws_ALC.Range("H55").Value = 932.6
$ws_ALC.Range("I55").Value = 1709
$ws_ALC.Range("J55").Value = 599.8570999999999
$ws_ALC.Range("K55").Value = 1709
$ws_ALC.Range("L55").Value = 599.8570999999999
$ws_ALC.Range("M55").Value = -1495
$ws_ALC.Range("N55").Value = -1027.8571

# ARM row 74
$ws_ARM.Range("H74").Value = 59954.15
$ws_ARM.Range("I74").Value = 67711.37
$ws_ARM.Range("J74").Value = 1775
$ws_ARM.Range("K74").Value = 67711.37
$ws_ARM.Range("L74").Value = 1775
$ws_ARM.Range("M74").Value = -66837.37
$ws_ARM.Range("N74").Value = -3523

# ARM row 77
$ws_ARM.Range("H77").Value = 59954.15
$ws_ARM.Range("I77").Value = 67711.37
$ws_ARM.Range("J77").Value = 1775
$ws_ARM.Range("K77").Value = 338556.85
$ws_ARM.Range("L77").Value = 8875
$ws_ARM.Range("M77").Value = -334188.85
$ws_ARM.Range("N77").Value = -17611

# CRP row 4
$ws_CRP.Range("H4").Value = 70002
$ws_CRP.Range("J4").Value = 70002
$ws_CRP.Range("L4").Value = 70002
$ws_CRP.Range("N4").Value = -70226

# CRP row 6
$ws_CRP.Range("H6").Value = 24237
$ws_CRP.Range("I6").Value = 0
$ws_CRP.Range("J6").Value = 24237
$ws_CRP.Range("K6").Value = 0
$ws_CRP.Range("L6").Value = 24237
$ws_CRP.Range("N6").Value = -24463
$ws_CRP.Range("M6").Value = ""

# CRP row 7
$ws_CRP.Range("H7").Value = 68.181816
$ws_CRP.Range("I7").Value = 71.666664
$ws_CRP.Range("J7").Value = 64
$ws_CRP.Range("K7").Value = 71.666664
$ws_CRP.Range("L7").Value = 64
$ws_CRP.Range("M7").Value = 41.333336
$ws_CRP.Range("N7").Value = -290

# CRP row 17
$ws_CRP.Range("H17").Value = 3887.5
$ws_CRP.Range("J17").Value = 4250
$ws_CRP.Range("L17").Value = 4250
$ws_CRP.Range("N17").Value = -4598

# CRP row 25
$ws_CRP.Range("H25").Value = 22363.916
$ws_CRP.Range("I25").Value = 1877.75
$ws_CRP.Range("J25").Value = 32607
$ws_CRP.Range("K25").Value = 1877.75
$ws_CRP.Range("L25").Value = 32607
$ws_CRP.Range("M25").Value = -1703.75
$ws_CRP.Range("N25").Value = -32955

# CRP row 31
$ws_CRP.Range("H31").Value = 18491.049
$ws_CRP.Range("I31").Value = 27886.703
$ws_CRP.Range("K31").Value = 27886.703
$ws_CRP.Range("M31").Value = -27591.703

# CRP row 34
$ws_CRP.Range("H34").Value = 18491.049
$ws_CRP.Range("I34").Value = 27886.703
$ws_CRP.Range("K34").Value = 27886.703
$ws_CRP.Range("M34").Value = -27684.703

# CRP row 41
$ws_CRP.Range("H41").Value = 8360.357
$ws_CRP.Range("I41").Value = 3256
$ws_CRP.Range("J41").Value = 11196.111
$ws_CRP.Range("K41").Value = 3256
$ws_CRP.Range("L41").Value = 11196.111
$ws_CRP.Range("M41").Value = -2828
$ws_CRP.Range("N41").Value = -12052.111

# CRP row 51
$ws_CRP.Range("H51").Value = 1000
$ws_CRP.Range("I51").Value = 1000
$ws_CRP.Range("K51").Value = 1000
$ws_CRP.Range("M51").Value = -264

# CRP row 59
$ws_CRP.Range("H59").Value = 11624
$ws_CRP.Range("I59").Value = 5000
$ws_CRP.Range("J59").Value = 12097.143
$ws_CRP.Range("K59").Value = 5000
$ws_CRP.Range("L59").Value = 12097.143
$ws_CRP.Range("M59").Value = -3855
$ws_CRP.Range("N59").Value = -14387.143

# CRP row 60
$ws_CRP.Range("H60").Value = 10202.667
$ws_CRP.Range("I60").Value = 3000
$ws_CRP.Range("K60").Value = 3000
$ws_CRP.Range("M60").Value = -2489

# CRP row 61
$ws_CRP.Range("H61").Value = 1000
$ws_CRP.Range("I61").Value = 1000
$ws_CRP.Range("K61").Value = 1000
$ws_CRP.Range("M61").Value = -652

# CRP row 68
$ws_CRP.Range("H68").Value = 27454.285
$ws_CRP.Range("J68").Value = 31363.334
$ws_CRP.Range("L68").Value = 31363.334
$ws_CRP.Range("N68").Value = -32861.334

# CRP row 71
$ws_CRP.Range("H71").Value = 27454.285
$ws_CRP.Range("J71").Value = 31363.334
$ws_CRP.Range("L71").Value = 94090.00199999999
$ws_CRP.Range("N71").Value = -101578.002

# CRP row 74
$ws_CRP.Range("H74").Value = 11890.4
$ws_CRP.Range("J74").Value = 11890.4
$ws_CRP.Range("L74").Value = 11890.4
$ws_CRP.Range("N74").Value = -13638.4

# CRP row 77
$ws_CRP.Range("H77").Value = 11890.4
$ws_CRP.Range("J77").Value = 11890.4
$ws_CRP.Range("L77").Value = 35671.2
$ws_CRP.Range("N77").Value = -44407.2

# CRP row 134
$ws_CRP.Range("H134").Value = 10870776
$ws_CRP.Range("I134").Value = 1122.2162
$ws_CRP.Range("K134").Value = 3366.6486
$ws_CRP.Range("M134").Value = -831.6486000000004

# CUL row 23
$ws_CUL.Range("H23").Value = 201.625
$ws_CUL.Range("I23").Value = 75
$ws_CUL.Range("J23").Value = 243.83333
$ws_CUL.Range("K23").Value = 225
$ws_CUL.Range("L23").Value = 731.49999
$ws_CUL.Range("M23").Value = 10
$ws_CUL.Range("N23").Value = -1201.49999

# CUL row 38
$ws_CUL.Range("H38").Value = 110.5
$ws_CUL.Range("I38").Value = 100.71429
$ws_CUL.Range("K38").Value = 302.14287
$ws_CUL.Range("M38").Value = 44.85712999999998

# CUL row 113
$ws_CUL.Range("H113").Value = 624.5
$ws_CUL.Range("I113").Value = 503
$ws_CUL.Range("J113").Value = 635.5454999999999
$ws_CUL.Range("K113").Value = 1509
$ws_CUL.Range("L113").Value = 1906.6365
$ws_CUL.Range("M113").Value = 661
$ws_CUL.Range("N113").Value = -6246.6365

# CUL row 122
$ws_CUL.Range("H122").Value = 376.64285
$ws_CUL.Range("I122").Value = 288.0909
$ws_CUL.Range("J122").Value = 701.3333
$ws_CUL.Range("K122").Value = 2592.8181
$ws_CUL.Range("L122").Value = 6311.9997
$ws_CUL.Range("M122").Value = -142.8181
$ws_CUL.Range("N122").Value = -11211.9997

# GSM row 5
$ws_GSM.Range("H5").Value = 13725.25
$ws_GSM.Range("J5").Value = 17967
$ws_GSM.Range("L5").Value = 17967
$ws_GSM.Range("N5").Value = -18191

# GSM row 112
$ws_GSM.Range("H112").Value = 86396.5
$ws_GSM.Range("J112").Value = 86396.5
$ws_GSM.Range("L112").Value = 86396.5
$ws_GSM.Range("N112").Value = -88612.5

# LTW row 46
$ws_LTW.Range("H46").Value = 1245.5
$ws_LTW.Range("I46").Value = 991
$ws_LTW.Range("J46").Value = 1500
$ws_LTW.Range("K46").Value = 991
$ws_LTW.Range("L46").Value = 1500
$ws_LTW.Range("M46").Value = -803
$ws_LTW.Range("N46").Value = -1876

# LTW row 55
$ws_LTW.Range("H55").Value = 478.1154
$ws_LTW.Range("I55").Value = 332.95
$ws_LTW.Range("J55").Value = 962
$ws_LTW.Range("K55").Value = 332.95
$ws_LTW.Range("L55").Value = 962
$ws_LTW.Range("M55").Value = -159.95
$ws_LTW.Range("N55").Value = -1308
